# Update the "想去人数" (want-to-go count) figures in column F
# for both the "展览" sheet and the "全部类型" sheet, which hold
# duplicate data.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 2199
    "F5"  = 13068
    "F10" = 1180
    "F11" = 978
    "F12" = 13756
    "F13" = 14331
    "F22" = 1089
    "F25" = 5389
    "F26" = 935
    "F28" = 306
    "F29" = 16
    "F30" = 28
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
